# "Double Bill Issue Fixed"
# - Stocks sheet: quantities reduced to reflect stock consumed by the new bills
# - Bills sheet: normalize Bal_Amt text "0.0" -> "0" on two existing rows, and
#   append 8 new bill rows (149-156) for three new transactions.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $text)
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- Stocks sheet (quantity column D) ----
$wsStocks = $wb.Worksheets.Item("Stocks")
$wsStocks.Range("D2").Value = 139.0
$wsStocks.Range("D3").Value = 402.0
$wsStocks.Range("D6").Value = 99.0
$wsStocks.Range("D7").Value = 199.0

# ---- Bills sheet ----
$wsBills = $wb.Worksheets.Item("Bills")

# Normalize existing Bal_Amt text from "0.0" to "0"
Set-TextCell $wsBills 144 5 "0"
Set-TextCell $wsBills 148 5 "0"

# New rows for three bills placed on 15-Dec-2020
$newRows = @(
    @(149, @("15-Dec-2020 09:58", "KKK",    "688", "775",  "0",   "XX1512148", "Stock Name One(5)")),
    @(150, @("15-Dec-2020 09:58", "KKK",    "688", "15",   "0",   "XX1512149", "New Liz(1)")),
    @(151, @("15-Dec-2020 09:58", "KKK",    "688", "690",  "0.0", "XX1512149", "New Liz(1),Bill Clearance 15Dec2020(1)")),
    @(152, @("15-Dec-2020 10:06", "Aakash", "364", "4600", "0",   "XX1512150", "Stock Name two(23)")),
    @(153, @("15-Dec-2020 10:06", "Aakash", "364", "25",   "0",   "XX1512151", "NewStock_Img(1)")),
    @(154, @("15-Dec-2020 10:10", "Aakash", "364", "3410", "0",   "XX1512151", "Stock Name One(22)")),
    @(155, @("15-Dec-2020 10:10", "Aakash", "364", "150",  "0",   "XX1512152", "Newss Sss(1)")),
    @(156, @("15-Dec-2020 10:10", "Aakash", "364", "3460", "0.0", "XX1512152", "Newss Sss(1),Bill Clearance 15Dec2020(1)"))
)

foreach ($entry in $newRows) {
    $rowNum = $entry[0]
    $values = $entry[1]
    for ($c = 1; $c -le 7; $c++) {
        Set-TextCell $wsBills $rowNum $c $values[$c - 1]
    }
}
